$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 8 ("Types of Variables" / byte,int,long,float,double,char lesson)
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)

# --- Shape 2: Title "Types of Variables:" -> "Types of Variables (according with C)"
$titleShape = $slide8.Shapes.Item(2)
$titleTr = $titleShape.TextFrame.TextRange
$titlePara1 = $titleTr.Paragraphs(1, 1)

# The trailing ":" run becomes a plain space (keeps its original 5400 size/formatting).
$lastChar = $titlePara1.Characters($titlePara1.Length, 1)
$lastChar.Text = " "

# Append the new "(according with C)" text as separate runs at 32pt, matching
# how the author typed them (parenthesis, word, space, word, space+paren).
$cur = $lastChar
$titlePieces = @("(", "according", " ", "with", " C)")
foreach ($piece in $titlePieces) {
    $null = $cur.InsertAfter($piece)
    $startIdx = $cur.Start + $cur.Length
    $cur = $titleTr.Characters($startIdx, $piece.Length)
    $cur.Font.Size = 32
}

# --- Shape 3: TextBox with the byte/int/long/float/double/char bullet list.
$bodyShape = $slide8.Shapes.Item(3)
$bodyTr = $bodyShape.TextFrame.TextRange

# First bullet ends with "...  8 bits ". Insert a brand-new bullet paragraph
# right after it: "short – A small number without floating-point".
$bytePara = $bodyTr.Paragraphs(1, 1)
$null = $bytePara.InsertAfter([char]13 + "short")

$shortStart = $bytePara.Start + $bytePara.Length
$shortRange = $bodyTr.Characters($shortStart, 5)
$shortRange.Font.Bold = $true

$cur = $shortRange
$shortPieces = @(" ", [char]8211 + " A ", "small", " ", "number", " ", "without", " ", "floating-point")
foreach ($piece in $shortPieces) {
    $null = $cur.InsertAfter($piece)
    $startIdx = $cur.Start + $cur.Length
    $cur = $bodyTr.Characters($startIdx, $piece.Length)
    $cur.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# Slide 9 ("Programming can be only made with 1's or 0's?")
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$questionShape = $slide9.Shapes.Item(3)
$questionTr = $questionShape.TextFrame.TextRange
$questionPara1 = $questionTr.Paragraphs(1, 1)

# " 0's" and "?" used to be two separate runs with identical formatting;
# normalise them into a single run (no visible text change).
$tailStart = $questionPara1.Length - 4
$tailRange = $questionPara1.Characters($tailStart, 5)
$apos = [string][char]8217
$tailRange.Text = " 0" + $apos + "s?"
